$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing cell styles for Price column cells we will force to Text
# so values like "1.000" or "0.07096" are not auto-converted to numbers.
$origStyles = @{}
$origStyles["D2"] = $ws.Range("D2").Style
$origStyles["D3"] = $ws.Range("D3").Style
$origStyles["D4"] = $ws.Range("D4").Style
$origStyles["D5"] = $ws.Range("D5").Style
$origStyles["D6"] = $ws.Range("D6").Style
$origStyles["D7"] = $ws.Range("D7").Style
$origStyles["D8"] = $ws.Range("D8").Style
$origStyles["D9"] = $ws.Range("D9").Style
$origStyles["D10"] = $ws.Range("D10").Style
$origStyles["D11"] = $ws.Range("D11").Style
$origStyles["D12"] = $ws.Range("D12").Style
$origStyles["D13"] = $ws.Range("D13").Style
$origStyles["D14"] = $ws.Range("D14").Style
$origStyles["D15"] = $ws.Range("D15").Style
$origStyles["D16"] = $ws.Range("D16").Style
$origStyles["D17"] = $ws.Range("D17").Style
$origStyles["D18"] = $ws.Range("D18").Style
$origStyles["D19"] = $ws.Range("D19").Style
$origStyles["D20"] = $ws.Range("D20").Style
$origStyles["D21"] = $ws.Range("D21").Style
$origStyles["D22"] = $ws.Range("D22").Style
$origStyles["D23"] = $ws.Range("D23").Style
$origStyles["D24"] = $ws.Range("D24").Style
$origStyles["D25"] = $ws.Range("D25").Style
$origStyles["D26"] = $ws.Range("D26").Style
$origStyles["D28"] = $ws.Range("D28").Style
$origStyles["D29"] = $ws.Range("D29").Style
$origStyles["D31"] = $ws.Range("D31").Style
$origStyles["D32"] = $ws.Range("D32").Style
$origStyles["D34"] = $ws.Range("D34").Style
$origStyles["D35"] = $ws.Range("D35").Style
$origStyles["D36"] = $ws.Range("D36").Style
$origStyles["D38"] = $ws.Range("D38").Style
$origStyles["D39"] = $ws.Range("D39").Style
$origStyles["D41"] = $ws.Range("D41").Style
$origStyles["D42"] = $ws.Range("D42").Style
$origStyles["D43"] = $ws.Range("D43").Style
$origStyles["D44"] = $ws.Range("D44").Style
$origStyles["D45"] = $ws.Range("D45").Style
$origStyles["D46"] = $ws.Range("D46").Style
$origStyles["D47"] = $ws.Range("D47").Style
$origStyles["D48"] = $ws.Range("D48").Style
$origStyles["D49"] = $ws.Range("D49").Style
$origStyles["D50"] = $ws.Range("D50").Style
$origStyles["D51"] = $ws.Range("D51").Style

# Temporarily format Price cells as Text so the literal strings are preserved
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update Price column values
$ws.Range("D2").Value = '30.310.24'
$ws.Range("D3").Value = '1.928.54'
$ws.Range("D4").Value = '1.000'
$ws.Range("D5").Value = '249.28'
$ws.Range("D6").Value = '0.7179'
$ws.Range("D7").Value = '1.000'
$ws.Range("D8").Value = '0.3206'
$ws.Range("D9").Value = '27.69'
$ws.Range("D10").Value = '0.07096'
$ws.Range("D11").Value = '0.7891'
$ws.Range("D12").Value = '0.07985'
$ws.Range("D13").Value = '1.928.89'
$ws.Range("D14").Value = '5.382'
$ws.Range("D15").Value = '94.82'
$ws.Range("D16").Value = '14.67'
$ws.Range("D17").Value = '30.309.93'
$ws.Range("D18").Value = '258.30'
$ws.Range("D19").Value = '0.000008101'
$ws.Range("D20").Value = '5.759'
$ws.Range("D21").Value = '2.185.53'
$ws.Range("D22").Value = '1.000'
$ws.Range("D23").Value = '1.001'
$ws.Range("D24").Value = '6.833'
$ws.Range("D25").Value = '9.537'
$ws.Range("D26").Value = '164.80'
$ws.Range("D28").Value = '2.273'
$ws.Range("D29").Value = '0.1269'
$ws.Range("D31").Value = '1.532'
$ws.Range("D32").Value = '4.401'
$ws.Range("D34").Value = '0.05146'
$ws.Range("D35").Value = '1.266'
$ws.Range("D36").Value = '0.7446'
$ws.Range("D38").Value = '0.01973'
$ws.Range("D39").Value = '2.799'
$ws.Range("D41").Value = '6.373'
$ws.Range("D42").Value = '0.4510'
$ws.Range("D43").Value = '1.996'
$ws.Range("D44").Value = '0.8481'
$ws.Range("D45").Value = '0.9999'
$ws.Range("D46").Value = '9.848'
$ws.Range("D47").Value = '100.59'
$ws.Range("D48").Value = '7.443'
$ws.Range("D49").Value = '36.82'
$ws.Range("D50").Value = '949.56'
$ws.Range("D51").Value = '0.4217'

# Restore original cell styles on the Price column cells
$ws.Range("D2").Style = $origStyles["D2"]
$ws.Range("D3").Style = $origStyles["D3"]
$ws.Range("D4").Style = $origStyles["D4"]
$ws.Range("D5").Style = $origStyles["D5"]
$ws.Range("D6").Style = $origStyles["D6"]
$ws.Range("D7").Style = $origStyles["D7"]
$ws.Range("D8").Style = $origStyles["D8"]
$ws.Range("D9").Style = $origStyles["D9"]
$ws.Range("D10").Style = $origStyles["D10"]
$ws.Range("D11").Style = $origStyles["D11"]
$ws.Range("D12").Style = $origStyles["D12"]
$ws.Range("D13").Style = $origStyles["D13"]
$ws.Range("D14").Style = $origStyles["D14"]
$ws.Range("D15").Style = $origStyles["D15"]
$ws.Range("D16").Style = $origStyles["D16"]
$ws.Range("D17").Style = $origStyles["D17"]
$ws.Range("D18").Style = $origStyles["D18"]
$ws.Range("D19").Style = $origStyles["D19"]
$ws.Range("D20").Style = $origStyles["D20"]
$ws.Range("D21").Style = $origStyles["D21"]
$ws.Range("D22").Style = $origStyles["D22"]
$ws.Range("D23").Style = $origStyles["D23"]
$ws.Range("D24").Style = $origStyles["D24"]
$ws.Range("D25").Style = $origStyles["D25"]
$ws.Range("D26").Style = $origStyles["D26"]
$ws.Range("D28").Style = $origStyles["D28"]
$ws.Range("D29").Style = $origStyles["D29"]
$ws.Range("D31").Style = $origStyles["D31"]
$ws.Range("D32").Style = $origStyles["D32"]
$ws.Range("D34").Style = $origStyles["D34"]
$ws.Range("D35").Style = $origStyles["D35"]
$ws.Range("D36").Style = $origStyles["D36"]
$ws.Range("D38").Style = $origStyles["D38"]
$ws.Range("D39").Style = $origStyles["D39"]
$ws.Range("D41").Style = $origStyles["D41"]
$ws.Range("D42").Style = $origStyles["D42"]
$ws.Range("D43").Style = $origStyles["D43"]
$ws.Range("D44").Style = $origStyles["D44"]
$ws.Range("D45").Style = $origStyles["D45"]
$ws.Range("D46").Style = $origStyles["D46"]
$ws.Range("D47").Style = $origStyles["D47"]
$ws.Range("D48").Style = $origStyles["D48"]
$ws.Range("D49").Style = $origStyles["D49"]
$ws.Range("D50").Style = $origStyles["D50"]
$ws.Range("D51").Style = $origStyles["D51"]

# Update remaining Coin / Link / Volume(1h) column values
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -4.95%  '
$ws.Range("E9").Value = '  -3.95%  '
$ws.Range("E10").Value = '  -3.15%  '
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("E12").Value = '  -2.08%  '
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("E14").Value = '  -2.81%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  -1.87%  '
$ws.Range("E25").Value = '  -3.20%  '
$ws.Range("E26").Value = '  +2.99%  '
$ws.Range("E27").Value = '  -1.78%  '
$ws.Range("E28").Value = '  -6.53%  '
$ws.Range("E29").Value = '  -4.13%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("E31").Value = '  -2.16%  '
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("E34").Value = '  -1.58%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("E37").Value = '  +1.03%  '
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("E40").Value = '  -3.11%  '
$ws.Range("E41").Value = '  -4.32%  '
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  -1.44%  '
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").Value = '  -4.03%  '
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("E50").Value = '  +8.15%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("E51").Value = '  +0.51%  '
